# Actualiza base de datos EC: intercambia los valores de "Valor Mora"
# entre la fila del periodo 2205 (fila 16) y la fila del periodo 2110 (fila 23).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("F16").Value = 36341
$ws.Range("F23").Value = 27861
